# Apply the changes described by the diff:
#  1. Rename the worksheet/workbook sheet entry from "o554F-HW20.xpc" to "o554F"
#  2. Tweak a few already-present values in row 13 (tiny precision corrections)
#  3. Append a new data row (row 16) with HKL index 14 / "HexGrid-60degTilt5degRes"
#     and its associated averaged-intensity values, extending the used range to A1:M16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "o554F"

# 2. Correct a handful of row-13 values (last-digit precision fixes)
$ws.Range("D13").Value = 0.9945358124397078
$ws.Range("J13").Value = 0.9945358124397078
$ws.Range("K13").Value = 0.994250140105256
$ws.Range("L13").Value = 0.9946371154641712

# 3. Add the new row 16.
# First copy the formatting of row 15's A-column cell (bold, bordered, centered style)
# so the new index cell A16 keeps the same look as A2:A15.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.166259248206962
$ws.Range("D16").Value = 0.6150372372946338
$ws.Range("E16").Value = 1.04570840306212
$ws.Range("F16").Value = 1.166259248206962
$ws.Range("G16").Value = 0.7983980656971881
$ws.Range("H16").Value = 1.126180299928017
$ws.Range("I16").Value = 1.086610671014488
$ws.Range("J16").Value = 0.6150372372946338
$ws.Range("K16").Value = 0.8303728201783769
$ws.Range("L16").Value = 0.9983160341926693
$ws.Range("M16").Value = 0.9730323208672349
